# Update countries & provincias Spain
# - Re-order "Republica de Macedonia" / "Senegal" in the country list
#   (Macedonia now appears right after Malasia, Senegal moves down one slot,
#   picking up fresh case numbers while the old Senegal numbers shift to
#   the row now labelled Senegal).
# - Refresh the "Datos actualizados" timestamp string.
# - Refresh COVID case counters for several countries/rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap order of Republica de Macedonia / Senegal (rows 80 & 81) -------
# Row 80 keeps the "Republica de Macedonia" label but gets fresh stats;
# Row 81 is now labelled "Senegal" and carries the stats that used to sit
# in row 80 before the refresh.
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("A81").Value = "Senegal"

# --- Footer timestamp ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 17:37"

# --- Row 4 : Estados Unidos -------------------------------------------------
$ws.Range("B4").Value = 3487635
$ws.Range("C4").Value = 8152
$ws.Range("E4").Value = 1798740
$ws.Range("G4").Value = 177
$ws.Range("H4").Value = 138424

# --- Row 6 : India -----------------------------------------------------------
$ws.Range("B6").Value = 933450
$ws.Range("C6").Value = 25805
$ws.Range("D6").Value = 590219
$ws.Range("E6").Value = 318950
$ws.Range("G6").Value = 554
$ws.Range("H6").Value = 24281

# --- Row 9 : Chile -----------------------------------------------------------
$ws.Range("B9").Value = 319493
$ws.Range("C9").Value = 1836
$ws.Range("D9").Value = 289220
$ws.Range("E9").Value = 23204
$ws.Range("G9").Value = 45
$ws.Range("H9").Value = 7069

# --- Row 43 --------------------------------------------------------------
$ws.Range("B43").Value = 47051
$ws.Range("C43").Value = 233
$ws.Range("D43").Value = 31550
$ws.Range("E43").Value = 13833
$ws.Range("G43").Value = 6
$ws.Range("H43").Value = 1668

# --- Row 44 --------------------------------------------------------------
$ws.Range("D44").Value = 42737
$ws.Range("E44").Value = 3866
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 27

# --- Row 45 --------------------------------------------------------------
$ws.Range("B45").Value = 46305
$ws.Range("C45").Value = 799
$ws.Range("D45").Value = 23134
$ws.Range("E45").Value = 22261
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 910

# --- Row 80 : Republica de Macedonia (new stats) ----------------------------
$ws.Range("B80").Value = 8332
$ws.Range("C80").Value = 135
$ws.Range("D80").Value = 4468
$ws.Range("E80").Value = 3475
$ws.Range("G80").Value = 4
$ws.Range("H80").Value = 389

# --- Row 81 : Senegal (carries the old row-80 stats) -------------------------
$ws.Range("B81").Value = 8243
$ws.Range("C81").Value = 45
$ws.Range("D81").Value = 5580
$ws.Range("E81").Value = 2513
$ws.Range("H81").Value = 150

# --- Row 88 --------------------------------------------------------------
$ws.Range("E88").Value = 5636
$ws.Range("G88").Value = 5
$ws.Range("H88").Value = 44

# --- Row 135 -------------------------------------------------------------
$ws.Range("B135").Value = 1268
$ws.Range("C135").Value = 49
$ws.Range("D135").Value = 373
$ws.Range("E135").Value = 886

# --- Row 143 -------------------------------------------------------------
$ws.Range("B143").Value = 1023
$ws.Range("C143").Value = 1
$ws.Range("E143").Value = 165

# --- Row 164 -------------------------------------------------------------
$ws.Range("B164").Value = 337
$ws.Range("C164").Value = 1
$ws.Range("E164").Value = 70
